$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B so the existing "dbExcel"/"WebExcel"
# columns shift right (B->C, C->D), making room for the new StatQuery column.
$ws.Columns.Item(2).Insert()

# The new column should be as wide as column A (matching the target layout,
# where columns A and B share the same 75.81640625 width). Columns C and D
# (previously B and C) keep whatever width they already had, so their exact
# widths are preserved untouched.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Header for the new column.
$ws.Range("B1").Value = "StatQuery"

# New stat-bar query text. Column B inherits the wrap-text style from column
# A (the newly-inserted column picks up the formatting of the column to its
# left), matching the target's s="1" on B2; set WrapText explicitly too so
# this does not depend on that inherited formatting.
$ws.Range("B2").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Bullmastiff']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$ws.Range("B2").WrapText = $true

# Update selection / view to match the new target state.
$ws.Range("A2").Select()
